$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.344.70'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -2.87%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.524.96'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -3.90%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '607.87'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.06'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.90%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.524.25'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.87%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -1.94%  '
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('E11').Value = '  -2.69%  '
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range('E13').Value = '  -4.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.122.43'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '31.86'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.526.70'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.78%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.378.29'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.37'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('E20').Value = '  -3.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '451.91'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.36'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.640'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '78.64'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.674.53'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.65%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('E28').Value = '  -3.81%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.31'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -7.30%  '
$ws.Range('E30').Value = '  -2.46%  '
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '25.93'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('E34').Value = '  -4.90%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.20'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.83%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.158'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.527.87'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.50%  '
$ws.Range('E38').Value = '  -5.36%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '176.65'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.60'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.94%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0877'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.894'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.23%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '29.23'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +8.65%  '
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.22'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('E51').Value = '  -3.74%  '